$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the row-index column (A) for rows 4..18 (rows 2 and 3 already have 1 and 2)
for ($row = 4; $row -le 18; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 1
}

# Update the account number in column F for rows 2..17 (row 18 keeps its own value)
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 6).Value = 5069929970
}

# Clear the FechaInicio value in N3 (removes the now-unused "07/04/2022" shared string)
$ws.Range("N3").Value = $null

# Update the sheet view: select F3:F17 (also clears the stale topLeftCell scroll position)
$ws.Range("F3:F17").Select() | Out-Null
